$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 346, shifting rows 346:378 down to 347:379
$ws.Rows.Item(346).Insert()

# Populate the newly inserted row 346 with the new price record
$ws.Cells.Item(346, 1).Value = 4
$ws.Cells.Item(346, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(346, 3).Value = "Los Lagos"
$ws.Cells.Item(346, 4).Value = 44578
$ws.Cells.Item(346, 5).Value = 10
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100102
$ws.Cells.Item(346, 8).Value = "Cítricos"
$ws.Cells.Item(346, 9).Value = 100102003
$ws.Cells.Item(346, 10).Value = "Limón"
$ws.Cells.Item(346, 11).Value = "Sin especificar"
$ws.Cells.Item(346, 12).Value = "1a plateado"
$ws.Cells.Item(346, 13).Value = 800
$ws.Cells.Item(346, 14).Value = 21000
$ws.Cells.Item(346, 15).Value = 22000
$ws.Cells.Item(346, 16).Value = 21500
$ws.Cells.Item(346, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(346, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(346, 19).Value = 1194
$ws.Cells.Item(346, 20).Value = 18
